$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that are stored as TEXT in the
# workbook. Force the text number format on the cells we touch so Excel
# keeps them as text instead of silently converting them to numbers.
$ws.Range("D2:D49").NumberFormat = "@"

# Simple price updates (column D only)
$ws.Range("D2").Value  = "248.96"
$ws.Range("D3").Value  = "22.48"
$ws.Range("D4").Value  = "5.403"
$ws.Range("D5").Value  = "0.05710"
$ws.Range("D6").Value  = "3.414"
$ws.Range("D7").Value  = "6.334"
$ws.Range("D8").Value  = "0.8138"
$ws.Range("D9").Value  = "0.9256"
$ws.Range("D10").Value = "0.1424"
$ws.Range("D11").Value = "0.07447"
$ws.Range("D14").Value = "0.09359"
$ws.Range("D15").Value = "3.738"
$ws.Range("D16").Value = "0.001597"
$ws.Range("D19").Value = "0.0005797"
$ws.Range("D20").Value = "0.006449"
$ws.Range("D21").Value = "0.005001"
$ws.Range("D22").Value = "0.001026"
$ws.Range("D25").Value = "2.166"
$ws.Range("D40").Value = "0.03994"

# Rows 41-43 got reshuffled: the three coins rotated position
# (old 41=KickToken, 42=BKEXToken, 43=CEJI -> new 41=BKEXToken, 42=CEJI, 43=KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1067"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002712"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.002983"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "0.007526"
$ws.Range("D45").Value = "0.00005898"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D49").Value = "0.00002103"
